$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws1.Range("H6").Value = 200629.8
$ws1.Range("I6").Value = 200629.8
$ws1.Range("K6").Value = 601889.3999999999
$ws1.Range("M6").Value = -601777.3999999999
$ws1.Range("H31").Value = 2006
$ws1.Range("I31").Value = 2006
$ws1.Range("K31").Value = 6018
$ws1.Range("M31").Value = -5788
$ws1.Range("H33").Value = 3965890.2
$ws1.Range("I33").Value = 5451617
$ws1.Range("J33").Value = 3953.1667
$ws1.Range("K33").Value = 5451617
$ws1.Range("L33").Value = 3953.1667
$ws1.Range("M33").Value = -5451388
$ws1.Range("N33").Value = -4411.1667
$ws1.Range("H63").Value = 28000
$ws1.Range("J63").Value = 28000
$ws1.Range("L63").Value = 28000
$ws1.Range("N63").Value = -29248
$ws1.Range("H66").Value = 28000
$ws1.Range("J66").Value = 28000
$ws1.Range("L66").Value = 84000
$ws1.Range("N66").Value = -90240
$ws1.Range("H70").Value = 1926.8
$ws1.Range("J70").Value = 2971
$ws1.Range("L70").Value = 8913
$ws1.Range("N70").Value = -9453
$ws1.Range("H73").Value = 1926.8
$ws1.Range("J73").Value = 2971
$ws1.Range("L73").Value = 8913
$ws1.Range("N73").Value = -10785
$ws1.Range("H76").Value = 3556.75
$ws1.Range("I76").Value = 3464
$ws1.Range("J76").Value = 3649.5
$ws1.Range("K76").Value = 3464
$ws1.Range("L76").Value = 3649.5
$ws1.Range("M76").Value = -3149
$ws1.Range("N76").Value = -4279.5
$ws1.Range("H79").Value = 3556.75
$ws1.Range("I79").Value = 3464
$ws1.Range("J79").Value = 3649.5
$ws1.Range("K79").Value = 3464
$ws1.Range("L79").Value = 3649.5
$ws1.Range("M79").Value = -2372
$ws1.Range("N79").Value = -5833.5
$ws1.Range("H107").Value = 14494457
$ws1.Range("I107").Value = 1298.3334
$ws1.Range("K107").Value = 1298.3334
$ws1.Range("M107").Value = 621.6666
$ws1.Range("H137").Value = 2089525.6
$ws1.Range("I137").Value = 3867.639
$ws1.Range("K137").Value = 11602.917
$ws1.Range("M137").Value = -9052.917000000001
$ws1.Range("H139").Value = 0
$ws1.Range("J139").Value = 0
$ws1.Range("L139").Value = 0
$ws1.Range("N139").ClearContents()

$ws2 = $wb.Worksheets.Item("ARM")
$ws2.Range("H45").Value = 79083.16
$ws2.Range("I45").Value = 101930.6
$ws2.Range("J45").Value = 2925
$ws2.Range("K45").Value = 101930.6
$ws2.Range("L45").Value = 2925
$ws2.Range("M45").Value = -101553.6
$ws2.Range("N45").Value = -3679
$ws2.Range("H52").Value = 29998
$ws2.Range("I52").Value = 29994
$ws2.Range("J52").Value = 30000
$ws2.Range("K52").Value = 29994
$ws2.Range("L52").Value = 30000
$ws2.Range("M52").Value = -29676
$ws2.Range("N52").Value = -30636
$ws2.Range("H61").Value = 1018936.7
$ws2.Range("I61").Value = 27524.568
$ws2.Range("J61").Value = 3745320
$ws2.Range("K61").Value = 27524.568
$ws2.Range("L61").Value = 3745320
$ws2.Range("M61").Value = -27312.568
$ws2.Range("N61").Value = -3745744
$ws2.Range("H74").Value = 403688.97
$ws2.Range("I74").Value = 2781.4
$ws2.Range("K74").Value = 2781.4
$ws2.Range("M74").Value = -1907.4
$ws2.Range("H77").Value = 403688.97
$ws2.Range("I77").Value = 2781.4
$ws2.Range("K77").Value = 13907
$ws2.Range("M77").Value = -9539
$ws2.Range("H105").Value = 81556.664
$ws2.Range("J105").Value = 81556.664
$ws2.Range("L105").Value = 81556.664
$ws2.Range("N105").Value = -88544.664
$ws2.Range("H122").Value = 2503.3333
$ws2.Range("I122").Value = 2004
$ws2.Range("K122").Value = 6012
$ws2.Range("M122").Value = -3562
$ws2.Range("H132").Value = 1504.091
$ws2.Range("I132").Value = 1042.2709
$ws2.Range("J132").Value = 4670.857
$ws2.Range("K132").Value = 3126.8127
$ws2.Range("L132").Value = 14012.571
$ws2.Range("M132").Value = -596.8126999999999
$ws2.Range("N132").Value = -19072.571
$ws2.Range("H136").Value = 1018936.7
$ws2.Range("I136").Value = 27524.568
$ws2.Range("J136").Value = 3745320
$ws2.Range("K136").Value = 82573.704
$ws2.Range("L136").Value = 11235960
$ws2.Range("M136").Value = -80023.704
$ws2.Range("N136").Value = -11241060

$ws3 = $wb.Worksheets.Item("BSM")
$ws3.Range("H134").Value = 21430284
$ws3.Range("I134").Value = 1385.7693
$ws3.Range("J134").Value = 56252244
$ws3.Range("K134").Value = 4157.3079
$ws3.Range("L134").Value = 168756732
$ws3.Range("M134").Value = -1622.3079
$ws3.Range("N134").Value = -168761802

$ws4 = $wb.Worksheets.Item("CRP")
$ws4.Range("H33").Value = 6727.2383
$ws4.Range("I33").Value = 1272.7273
$ws4.Range("J33").Value = 12727.2
$ws4.Range("K33").Value = 1272.7273
$ws4.Range("L33").Value = 12727.2
$ws4.Range("M33").Value = -893.7273
$ws4.Range("N33").Value = -13485.2
$ws4.Range("H42").Value = 13100
$ws4.Range("I42").Value = 10000
$ws4.Range("J42").Value = 16200
$ws4.Range("K42").Value = 10000
$ws4.Range("L42").Value = 16200
$ws4.Range("M42").Value = -9407
$ws4.Range("N42").Value = -17386
$ws4.Range("H58").Value = 941
$ws4.Range("I58").Value = 941
$ws4.Range("K58").Value = 941
$ws4.Range("M58").Value = -738
$ws4.Range("H132").Value = 2549.1
$ws4.Range("I132").Value = 1981
$ws4.Range("K132").Value = 5943
$ws4.Range("M132").Value = -3413
$ws4.Range("H134").Value = 2593.4
$ws4.Range("I134").Value = 2416.8333
$ws4.Range("K134").Value = 7250.499899999999
$ws4.Range("M134").Value = -4715.499899999999
$ws4.Range("H136").Value = 941
$ws4.Range("I136").Value = 941
$ws4.Range("K136").Value = 2823
$ws4.Range("M136").Value = -273

$ws5 = $wb.Worksheets.Item("CUL")
$ws5.Range("H4").Value = 2539216.5
$ws5.Range("I4").Value = 5000721.5
$ws5.Range("J4").Value = 734112.9
$ws5.Range("K4").Value = 15002164.5
$ws5.Range("L4").Value = 2202338.7
$ws5.Range("M4").Value = -15002052.5
$ws5.Range("N4").Value = -2202562.7
$ws5.Range("H5").Value = 1852.625
$ws5.Range("I5").Value = 1205.25
$ws5.Range("K5").Value = 3615.75
$ws5.Range("M5").Value = -3503.75
$ws5.Range("H7").Value = 160.90909
$ws5.Range("J7").Value = 125.25
$ws5.Range("L7").Value = 375.75
$ws5.Range("N7").Value = -599.75
$ws5.Range("H59").Value = 10000
$ws5.Range("J59").Value = 10000
$ws5.Range("L59").Value = 30000
$ws5.Range("N59").Value = -31080
$ws5.Range("H60").Value = 245.83333
$ws5.Range("I60").Value = 235
$ws5.Range("K60").Value = 705
$ws5.Range("M60").Value = -454
$ws5.Range("H123").Value = 3885
$ws5.Range("I123").Value = 3885
$ws5.Range("K123").Value = 11655
$ws5.Range("M123").Value = -9205
$ws5.Range("H135").Value = 1852.625
$ws5.Range("I135").Value = 1205.25
$ws5.Range("K135").Value = 10847.25
$ws5.Range("M135").Value = -8312.25
$ws5.Range("H138").Value = 4303.1577
$ws5.Range("I138").Value = 4138.1875
$ws5.Range("K138").Value = 12414.5625
$ws5.Range("M138").Value = -7274.5625

$ws6 = $wb.Worksheets.Item("GSM")
$ws6.Range("H5").Value = 13004.667
$ws6.Range("H80").Value = 31261408
$ws6.Range("I80").Value = 6750.6
$ws6.Range("K80").Value = 6750.6
$ws6.Range("M80").Value = -5752.6
$ws6.Range("H83").Value = 31261408
$ws6.Range("I83").Value = 6750.6
$ws6.Range("K83").Value = 33753
$ws6.Range("M83").Value = -28761

$ws7 = $wb.Worksheets.Item("LTW")
$ws7.Range("H12").Value = 0
$ws7.Range("I12").Value = 0
$ws7.Range("K12").Value = 0
$ws7.Range("M12").ClearContents()
$ws7.Range("H93").Value = 1361.6923
$ws7.Range("I93").Value = 1488
$ws7.Range("J93").Value = 667
$ws7.Range("K93").Value = 1488
$ws7.Range("L93").Value = 667
$ws7.Range("M93").Value = -240
$ws7.Range("N93").Value = -3163
$ws7.Range("H122").Value = 3702
$ws7.Range("I122").Value = 3124.889
$ws7.Range("K122").Value = 9374.667000000001
$ws7.Range("M122").Value = -6924.667000000001
$ws7.Range("H136").Value = 114730.336
$ws7.Range("J136").Value = 2595
$ws7.Range("L136").Value = 7785
$ws7.Range("N136").Value = -12885

$ws8 = $wb.Worksheets.Item("WVR")
$ws8.Range("H29").Value = 0
$ws8.Range("J29").Value = 0
$ws8.Range("L29").Value = 0
$ws8.Range("H81").Value = 6623327.5
$ws8.Range("I81").Value = 7560548.5
$ws8.Range("K81").Value = 15121097
$ws8.Range("M81").Value = -15120036
$ws8.Range("H84").Value = 6623327.5
$ws8.Range("I84").Value = 7560548.5
$ws8.Range("K84").Value = 75605485
$ws8.Range("M84").Value = -75600181
$ws8.Range("H107").Value = 2042859.2
$ws8.Range("I107").Value = 1633
$ws8.Range("K107").Value = 4899
$ws8.Range("M107").Value = -2979
$ws8.Range("H132").Value = 2463
$ws8.Range("I132").Value = 2104.7693
$ws8.Range("K132").Value = 6314.3079
$ws8.Range("M132").Value = -3784.3079
$ws8.Range("H136").Value = 1071.4445
$ws8.Range("I136").Value = 909
$ws8.Range("J136").Value = 1274.5
$ws8.Range("K136").Value = 2727
$ws8.Range("L136").Value = 3823.5
$ws8.Range("M136").Value = -177
$ws8.Range("N136").Value = -8923.5
$ws8.Range("N29").ClearContents()
